$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.713.97'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.278.85'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '309.45'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.13%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.01'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.598'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.34'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0895'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('E12').Value = '  -1.56%  '
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.97'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '2.623.91'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').Value = '2.273.83'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '42.604.85'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.22'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('E23').Value = '  -6.63%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '261.10'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.16'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.78%  '
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.58'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.32'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.83'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +13.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.02'
$ws.Range('D30').ClearFormats()
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '164.05'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '35.33'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0847'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('E36').Value = '  -4.43%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.44'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0345'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.61%  '
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.70'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.71%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.55'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '98.89'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +9.71%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.01'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '68.03'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '11.86'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').Value = '1.704.36'
$ws.Range('E47').Value = '  +6.90%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '109.32'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '76.14'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.58'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.84%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.10'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.07%  '
